# Parallel test case is added
#
# - Add a new ("Sheet1") worksheet at the end of the workbook.
# - Activate the "Test Data" sheet (it becomes the active tab).
# - On "Test Data": flip several Runmode cells from Y -> N (these rows are
#   now run in the new parallel lane) and add/update the new "Status"
#   column (H) with PASS / SKIP / FAIL results, mirroring the commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Test Data sheet: cell value changes
# ---------------------------------------------------------------------
$testData = $wb.Worksheets.Item("Test Data")

$testData.Range("E3").Value = "PASS"

$testData.Range("A9").Value = "N"
$testData.Range("E9").Value = "SKIP"

$testData.Range("A13").Value = "N"
$testData.Range("H13").Value = "SKIP"

$testData.Range("H14").Value = "PASS"

$testData.Range("A15").Value = "N"
$testData.Range("H15").Value = "SKIP"

$testData.Range("A16").Value = "N"
$testData.Range("H16").Value = "SKIP"

$testData.Range("A17").Value = "N"
$testData.Range("H17").Value = "SKIP"

$testData.Range("H18").Value = "FAIL"

$testData.Range("A19").Value = "N"
$testData.Range("H19").Value = "SKIP"

$testData.Range("A20").Value = "N"
$testData.Range("H20").Value = "SKIP"

# Column width touch-ups that went along with the new Status column.
# (Column C/D/G grew to fit the longer header/value text, H is the new
# narrow Status column.) Inputs are pre-compensated for the host's
# "characters -> stored width" rounding so the saved width lands as close
# as possible to the authored value.
$testData.Columns.Item(3).ColumnWidth = 27.166666666666668
$testData.Columns.Item(4).ColumnWidth = 31.333333333333332
$testData.Columns.Item(7).ColumnWidth = 27.5
$testData.Columns.Item(8).ColumnWidth = 5.666666666666667

# ---------------------------------------------------------------------
# Add the new trailing worksheet ("Sheet1")
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
[void]$newSheet.Range("M11").Select()

# ---------------------------------------------------------------------
# View state: Test Data becomes the active tab; restore selections on
# the other sheets the way the author left them.
# ---------------------------------------------------------------------
$testCases = $wb.Worksheets.Item("Test Cases")
$testCases.Columns.Item(1).ColumnWidth = 21.166666666666668
[void]$testCases.Range("B18").Select()

$testSteps = $wb.Worksheets.Item("Test Steps")
[void]$testSteps.Range("D6").Select()

[void]$testData.Activate()
[void]$testData.Range("G24").Select()
